# Update view-count-like figures (column F) across the four sheets to match
# the latest scrape snapshot (commit "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 2929
$ws1.Range("F10").Value = 6852
$ws1.Range("F12").Value = 56
$ws1.Range("F15").Value = 1484
$ws1.Range("F17").Value = 2224
$ws1.Range("F18").Value = 1470
$ws1.Range("F20").Value = 105
$ws1.Range("F21").Value = 1108
$ws1.Range("F23").Value = 173
$ws1.Range("F25").Value = 1693
$ws1.Range("F29").Value = 35
$ws1.Range("F30").Value = 1659
$ws1.Range("F31").Value = 1205
$ws1.Range("F36").Value = 418
$ws1.Range("F38").Value = 2450
$ws1.Range("F39").Value = 2707
$ws1.Range("F41").Value = 185
$ws1.Range("F49").Value = 413

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F23").Value = 465

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value  = 1669
$ws3.Range("F13").Value = 1416
$ws3.Range("F14").Value = 7298

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 2929
$ws4.Range("F7").Value  = 1669
$ws4.Range("F10").Value = 6852
$ws4.Range("F17").Value = 1416
$ws4.Range("F19").Value = 2224
$ws4.Range("F20").Value = 1470
$ws4.Range("F22").Value = 105
$ws4.Range("F23").Value = 1108
$ws4.Range("F27").Value = 1693
$ws4.Range("F30").Value = 35
$ws4.Range("F31").Value = 1659
$ws4.Range("F32").Value = 1205
$ws4.Range("F36").Value = 465
$ws4.Range("F37").Value = 418
$ws4.Range("F39").Value = 2450
$ws4.Range("F40").Value = 2707
$ws4.Range("F42").Value = 185
$ws4.Range("F48").Value = 413

$wb.Save()
